$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.310.26"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.869.33"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4697"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2868"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06581"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07862"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "1.870.47"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6920"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.101"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "267.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "30.293.53"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007671"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "2.116.46"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.243"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.204"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.384"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.948"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09888"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.355"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.358"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.463"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.070"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04755"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7023"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.725"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01875"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.806"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.210"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.954"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8439"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4175"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.125"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "938.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.115"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05683"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
